# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~18:38) to the
# SeniorConnect master log workbook across the affected sheets:
#   - Proximity      : new ENTER event at the Living Room Main Door
#   - Camera         : new "Image Captured" event at the Living Room Main Door
#   - mmWave(BR)     : new bedroom presence readings (numeric distance/value)
#   - mmWave(HR)     : new bedroom presence readings (numeric distance/value)
#   - mmWave(InBed)  : new bedroom in-bed / out-of-bed readings (text value)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Proximity: single new row (row 2)
# ---------------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$wsProximity.Cells.Item(2, 1).NumberFormat = "@"
$wsProximity.Cells.Item(2, 1).Value = "2026-01-28"
$wsProximity.Cells.Item(2, 2).Value = "18:38:33"
$wsProximity.Cells.Item(2, 3).Value = "18:00"
$wsProximity.Cells.Item(2, 4).Value = "Living Room Main Door"
$wsProximity.Cells.Item(2, 5).Value = "ENTER"
$wsProximity.Cells.Item(2, 6).Value = "User ENTERED Living Room Main Door"

# ---------------------------------------------------------------------------
# Camera: single new row (row 2)
# ---------------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")
$wsCamera.Cells.Item(2, 1).NumberFormat = "@"
$wsCamera.Cells.Item(2, 1).Value = "2026-01-28"
$wsCamera.Cells.Item(2, 2).Value = "18:38:35"
$wsCamera.Cells.Item(2, 3).Value = "18:00"
$wsCamera.Cells.Item(2, 4).Value = "Living Room Main Door"
$wsCamera.Cells.Item(2, 5).Value = "Image Captured"
$wsCamera.Cells.Item(2, 6).Value = "Active"

# ---------------------------------------------------------------------------
# mmWave(InBed): new rows 22-30 (Value column is text: In Bed / Out of Bed)
# ---------------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
$inBedData = @(
    @("2026-01-28", "18:38:05", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:06", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:07", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:08", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:09", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:10", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:11", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:15", "18:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-01-28", "18:38:34", "18:00", "Bedroom", "Out of Bed", "Empty")
)
$r = 22
foreach ($row in $inBedData) {
    $wsInBed.Cells.Item($r, 1).NumberFormat = "@"
    $wsInBed.Cells.Item($r, 1).Value = $row[0]
    $wsInBed.Cells.Item($r, 2).Value = $row[1]
    $wsInBed.Cells.Item($r, 3).Value = $row[2]
    $wsInBed.Cells.Item($r, 4).Value = $row[3]
    $wsInBed.Cells.Item($r, 5).Value = $row[4]
    $wsInBed.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# mmWave(BR): new rows 22-30 (Value column is numeric)
# ---------------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")
$brData = @(
    @("2026-01-28", "18:38:05", "18:00", "Bedroom", 0,  "Occupied"),
    @("2026-01-28", "18:38:07", "18:00", "Bedroom", 20, "Occupied"),
    @("2026-01-28", "18:38:08", "18:00", "Bedroom", 4,  "Occupied"),
    @("2026-01-28", "18:38:08", "18:00", "Bedroom", 17, "Occupied"),
    @("2026-01-28", "18:38:09", "18:00", "Bedroom", 2,  "Occupied"),
    @("2026-01-28", "18:38:11", "18:00", "Bedroom", 14, "Occupied"),
    @("2026-01-28", "18:38:11", "18:00", "Bedroom", 2,  "Occupied"),
    @("2026-01-28", "18:38:16", "18:00", "Bedroom", 1,  "Occupied"),
    @("2026-01-28", "18:38:34", "18:00", "Bedroom", 0,  "Empty")
)
$r = 22
foreach ($row in $brData) {
    $wsBR.Cells.Item($r, 1).NumberFormat = "@"
    $wsBR.Cells.Item($r, 1).Value = $row[0]
    $wsBR.Cells.Item($r, 2).Value = $row[1]
    $wsBR.Cells.Item($r, 3).Value = $row[2]
    $wsBR.Cells.Item($r, 4).Value = $row[3]
    $wsBR.Cells.Item($r, 5).Value = $row[4]
    $wsBR.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# mmWave(HR): new rows 22-30 (Value column is numeric)
# ---------------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")
$hrData = @(
    @("2026-01-28", "18:38:05", "18:00", "Bedroom", 0,  "Occupied"),
    @("2026-01-28", "18:38:07", "18:00", "Bedroom", 68, "Occupied"),
    @("2026-01-28", "18:38:08", "18:00", "Bedroom", 52, "Occupied"),
    @("2026-01-28", "18:38:08", "18:00", "Bedroom", 65, "Occupied"),
    @("2026-01-28", "18:38:09", "18:00", "Bedroom", 50, "Occupied"),
    @("2026-01-28", "18:38:11", "18:00", "Bedroom", 62, "Occupied"),
    @("2026-01-28", "18:38:11", "18:00", "Bedroom", 50, "Occupied"),
    @("2026-01-28", "18:38:15", "18:00", "Bedroom", 49, "Occupied"),
    @("2026-01-28", "18:38:34", "18:00", "Bedroom", 0,  "Empty")
)
$r = 22
foreach ($row in $hrData) {
    $wsHR.Cells.Item($r, 1).NumberFormat = "@"
    $wsHR.Cells.Item($r, 1).Value = $row[0]
    $wsHR.Cells.Item($r, 2).Value = $row[1]
    $wsHR.Cells.Item($r, 3).Value = $row[2]
    $wsHR.Cells.Item($r, 4).Value = $row[3]
    $wsHR.Cells.Item($r, 5).Value = $row[4]
    $wsHR.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
